# Daily auto push: insert a new reading for 2026/01/12 at the top of the
# "time-series tail" block (row 634), pushing every later row down by one.
# The last existing row (formerly 675, date 2027/01/05) ends up at 676 and
# the sheet's used range grows from D675 to D676.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 634..675 down to 635..676, leaving row 634 blank for the new entry.
$ws.Rows.Item(634).Insert()

# New row 634: 2026/01/12 (月), time 13, ranking 158.
# Force the date column to text so it stores the literal "2026/01/12" string
# (matching the rest of column A) instead of being auto-parsed into a date.
$ws.Range("A634").NumberFormat = "@"
$ws.Range("A634").Value = "2026/01/12"
$ws.Range("B634").Value = "月"
$ws.Range("C634").Value = 13
$ws.Range("D634").Value = 158
